$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Bump the date in A1 by one day (45310 -> 45311 serial date)
$ws.Range("A1").Value = 45311

# Update unit price (PRECIO C/U) column D for rows 14-21
$ws.Range("D14").Value = 38.5
$ws.Range("D15").Value = 50.1
$ws.Range("D16").Value = 60
$ws.Range("D17").Value = 95.7
$ws.Range("D18").Value = 119.4
$ws.Range("D19").Value = 149.1
$ws.Range("D20").Value = 167
$ws.Range("D21").Value = 185
